$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-34: update Price (D) and Volume(1h) (E) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.130.00"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.992.48"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.17"
$ws.Range("E5").Value = "  +12.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.80"
$ws.Range("E6").Value = "  +10.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.686"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.755"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.70"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.05"
$ws.Range("E13").Value = "  +4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.627.16"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.996.58"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.27"
$ws.Range("E16").Value = "  +8.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.55"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.808.57"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.34"
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.82"
$ws.Range("E22").Value = "  +13.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.97"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.46"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.42"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.34"
$ws.Range("E26").Value = "  +13.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.40"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.46"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.51"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.94"
$ws.Range("E32").Value = "  +4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.132"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "49.14"
$ws.Range("E34").Value = "  -3.73%  "

# Row 35 & 36: swap OKB/Bittensor content (with updated values)
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "70.99"
$ws.Range("E35").Value = "  +7.90%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "664.00"
$ws.Range("E36").Value = "  -2.09%  "

# Rows 37-38: update Price (D) and Volume(1h) (E) values
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0902"
$ws.Range("E37").Value = "  +10.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.439"
$ws.Range("E38").Value = "  -0.51%  "

# Row 39 & 40: swap Kaspa/ThetaToken content (with updated values)
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").Value = "  -0.78%  "

# Rows 41-51: update Price (D) and Volume(1h) (E) values
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.73"
$ws.Range("E45").Value = "  +4.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.151"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.39"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.910.04"
$ws.Range("E49").Value = "  +12.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.07"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  +5.06%  "

